# Apply updated descriptive-statistics values to the "Tbl2" worksheet
# (commit: "adding qualitative categories to main descriptive table").
# The underlying analysis was re-run after adding qualitative category
# variables, which shifted the correlation / mean statistics in the
# existing cells below. No new rows or columns are introduced here -
# only the numeric results already on the sheet are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tbl2")
$ws.Activate()

$ws.Range("C5").Value2 = 0.956730769230769
$ws.Range("D5").Value2 = 0.892742796157951
$ws.Range("E5").Value2 = 0.924719701014415
$ws.Range("C6").Value2 = 0.769230769230769
$ws.Range("D6").Value2 = 0.481323372465315
$ws.Range("E6").Value2 = 0.625200213561132
$ws.Range("F6").Value2 = 0.368506798064694
$ws.Range("C7").Value2 = 0.0134428418803419
$ws.Range("D7").Value2 = 0.0711980256136606
$ws.Range("E7").Value2 = 0.0423358515750134
$ws.Range("F7").Value2 = -0.303398339202107
$ws.Range("G7").Value2 = -0.351455865995926
$ws.Range("C9").Value2 = 0.981837606837607
$ws.Range("D9").Value2 = 0.9232
$ws.Range("E9").Value2 = 0.952495329597011
$ws.Range("C10").Value2 = 0.838675213675214
$ws.Range("D10").Value2 = 0.5664
$ws.Range("E10").Value2 = 0.70242860955431
$ws.Range("I10").Value2 = 0.343116999926229
$ws.Range("C11").Value2 = 0.00790598290598291
$ws.Range("D11").Value2 = 0.0249600533333333
$ws.Range("E11").Value2 = 0.0164398452095009
$ws.Range("I11").Value2 = -0.594126146367684
$ws.Range("J11").Value2 = -0.51623280717746
$ws.Range("D13").Value2 = 0.621131270010672
$ws.Range("E13").Value2 = 0.643886812600107
$ws.Range("F13").Value2 = 0.12376790047938
$ws.Range("G13").Value2 = 0.198105161441526
$ws.Range("H13").Value2 = -0.102311120162046
$ws.Range("I13").Value2 = 0.103837549767532
$ws.Range("J13").Value2 = 0.160972424171552
$ws.Range("K13").Value2 = -0.144229225401533
$ws.Range("D14").Value2 = 0.510138740661686
$ws.Range("E14").Value2 = 0.498665242925788
$ws.Range("F14").Value2 = 0.049
$ws.Range("G14").Value2 = 0.005
$ws.Range("H14").Value2 = 0.039
$ws.Range("I14").Value2 = 0.036
$ws.Range("J14").Value2 = 0.002
$ws.Range("K14").Value2 = -0.08
$ws.Range("D15").Value2 = 0.0921157024384299
$ws.Range("E15").Value2 = 0.0225364443884533
$ws.Range("F15").Value2 = -0.00329059557576967
$ws.Range("G15").Value2 = -0.0392272997449865
$ws.Range("H15").Value2 = -0.0480861604427693
$ws.Range("I15").Value2 = 0.0168068757679965
$ws.Range("J15").Value2 = -0.0166326814231581
$ws.Range("K15").Value2 = -0.00684092803989187
$ws.Range("D16").Value2 = 0.0070971184631805
$ws.Range("E16").Value2 = -0.0406567004805124
$ws.Range("F16").Value2 = -0.185968922377949
$ws.Range("G16").Value2 = -0.239321783172534
$ws.Range("H16").Value2 = -0.000435416962980169
$ws.Range("I16").Value2 = -0.218509776636481
$ws.Range("J16").Value2 = -0.230293303138601
$ws.Range("K16").Value2 = 0.289279450910045
$ws.Range("D17").Value2 = 2.25026680896478
$ws.Range("E17").Value2 = 2.17645488521089
$ws.Range("F17").Value2 = -0.0349595818002763
$ws.Range("G17").Value2 = -0.00527012637152524
$ws.Range("H17").Value2 = 0.0654573703959077
$ws.Range("I17").Value2 = -0.00925197681418996
$ws.Range("J17").Value2 = 0.0220595407837822
$ws.Range("K17").Value2 = -0.005154042648395
$ws.Range("C18").Value2 = 0.358974358974359
$ws.Range("D18").Value2 = 0.702241195304162
$ws.Range("E18").Value2 = 0.530699412706887
$ws.Range("F18").Value2 = 0.018
$ws.Range("G18").Value2 = -0.237
$ws.Range("H18").Value2 = -0.008
$ws.Range("I18").Value2 = 0.083
$ws.Range("J18").Value2 = -0.179
$ws.Range("K18").Value2 = -0.027
$ws.Range("E19").Value2 = 0.499733048585157
$ws.Range("F19").Value2 = 0.121261398612089
$ws.Range("G19").Value2 = 0.297381220063287
$ws.Range("H19").Value2 = -0.195510136601681
$ws.Range("I19").Value2 = 0.137830877260125
$ws.Range("J19").Value2 = 0.29777023906932
$ws.Range("K19").Value2 = -0.200825518271152

# Reflect the reviewer's final cursor position/selection on the sheet
# (was B15 before the refresh; ends on G9 after reviewing the updated
# "Curated Replication Results" verification correlations).
$ws.Range("G9").Select()
